$wb = $excel.ActiveWorkbook

# Rename "Klas3" to "Klas2" (note tracking is now kept for Klas2)
$ws2 = $wb.Worksheets.Item("Klas3")
$ws2.Name = "Klas2"

# Add the new "Nota" column (C) with the same header style as A1/B1
$ws2.Range("C1").Value = "Nota"
$ws2.Range("A1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("C1").Value = "Nota"

# Fill the new column with an (unchecked) boolean for each student row
for ($r = 2; $r -le 24; $r++) {
    $ws2.Cells.Item($r, 3).Value = $false
}

# Klas2 becomes the active sheet/tab, with E7 selected
$ws2.Activate()
$ws2.Range("E7").Select()
